$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain plain text so number-like values
# (e.g. "301.90", "1.000", "0.00001210") keep their exact formatting
# instead of being auto-converted to numeric values by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "23.237.70"
$ws.Range("E2").Value = "  -0.71%  "
$ws.Range("D3").Value = "1.612.43"
$ws.Range("E3").Value = "  -1.00%  "
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.31%  "
$ws.Range("D5").Value = "1.001"
$ws.Range("E5").Value = "  +0.16%  "
$ws.Range("D6").Value = "301.90"
$ws.Range("E6").Value = "  -0.78%  "
$ws.Range("D7").Value = "0.3781"
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("D8").Value = "51.79"
$ws.Range("E8").Value = "  -0.25%  "
$ws.Range("D9").Value = "0.3526"
$ws.Range("E9").Value = "  -2.92%  "
$ws.Range("D10").Value = "0.08074"
$ws.Range("E10").Value = "  -0.30%  "
$ws.Range("D11").Value = "1.202"
$ws.Range("E11").Value = "  -2.45%  "
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  +0.28%  "
$ws.Range("D13").Value = "22.05"
$ws.Range("E13").Value = "  -2.77%  "
$ws.Range("D14").Value = "6.347"
$ws.Range("E14").Value = "  -3.07%  "
$ws.Range("D15").Value = "7.250"
$ws.Range("E15").Value = "  +0.50%  "
$ws.Range("D16").Value = "0.00001210"
$ws.Range("E16").Value = "  -2.98%  "
$ws.Range("D17").Value = "1.589.34"
$ws.Range("E17").Value = "  -2.36%  "
$ws.Range("D18").Value = "93.73"
$ws.Range("E18").Value = "  +0.19%  "
$ws.Range("D19").Value = "0.06896"
$ws.Range("E19").Value = "  -0.07%  "
$ws.Range("D20").Value = "6.448"
$ws.Range("E20").Value = "  +0.51%  "
$ws.Range("D21").Value = "1.002"
$ws.Range("E21").Value = "  +0.20%  "
$ws.Range("D22").Value = "17.16"
$ws.Range("E22").Value = "  -4.12%  "
$ws.Range("D23").Value = "12.28"
$ws.Range("E23").Value = "  -3.53%  "
$ws.Range("D24").Value = "23.230.08"
$ws.Range("E24").Value = "  -0.76%  "
$ws.Range("D25").Value = "2.535"
$ws.Range("E25").Value = "  +3.61%  "
$ws.Range("D26").Value = "3.052"
$ws.Range("E26").Value = "  -6.28%  "
$ws.Range("D27").Value = "20.80"
$ws.Range("D28").Value = "151.09"
$ws.Range("E28").Value = "  +1.58%  "
$ws.Range("E29").Value = "  -0.81%  "
$ws.Range("D30").Value = "131.90"
$ws.Range("E30").Value = "  -1.86%  "
$ws.Range("D31").Value = "1.775.26"
$ws.Range("E31").Value = "  -1.90%  "
$ws.Range("D32").Value = "1.066"
$ws.Range("E32").Value = "  +11.95%  "
$ws.Range("D33").Value = "6.436"
$ws.Range("E33").Value = "  -5.54%  "
$ws.Range("D34").Value = "2.103"
$ws.Range("E34").Value = "  -9.22%  "
$ws.Range("D35").Value = "11.38"
$ws.Range("E35").Value = "  +2.84%  "
$ws.Range("D36").Value = "0.02701"
$ws.Range("E36").Value = "  -3.39%  "
$ws.Range("D37").Value = "0.08664"
$ws.Range("E37").Value = "  -1.78%  "
$ws.Range("D38").Value = "0.2448"
$ws.Range("D39").Value = "0.06902"
$ws.Range("E39").Value = "  -3.83%  "
$ws.Range("D40").Value = "5.836"
$ws.Range("E40").Value = "  -4.46%  "
$ws.Range("D41").Value = "1.317"
$ws.Range("E41").Value = "  -3.11%  "
$ws.Range("D42").Value = "0.6849"
$ws.Range("E42").Value = "  -3.09%  "
$ws.Range("D43").Value = "11.95"
$ws.Range("E43").Value = "  -2.99%  "
$ws.Range("D44").Value = "15.15"
$ws.Range("E44").Value = "  -6.40%  "
$ws.Range("D45").Value = "1.000"
$ws.Range("E45").Value = "  +0.14%  "
$ws.Range("D46").Value = "0.6290"
$ws.Range("E46").Value = "  -2.60%  "
$ws.Range("D47").Value = "3.934"
$ws.Range("E47").Value = "  -1.57%  "
$ws.Range("D48").Value = "2.241"
$ws.Range("E48").Value = "  -3.71%  "
$ws.Range("D49").Value = "0.07865"
$ws.Range("E49").Value = "  -1.60%  "
$ws.Range("D50").Value = "127.69"
$ws.Range("E50").Value = "  +1.53%  "
$ws.Range("D51").Value = "1.166"
$ws.Range("E51").Value = "  -3.11%  "
